$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "kora baki" update: insert two new columns into the asset schedule.
# New column F = "Cost of asset sold" (inserted before the old "Total" column),
# New column H = "Current Balance" (inserted before the old "Rate (Pre MAB Assets)" column).
# Column insertion shifts all existing F:P content (and formatting/types, e.g. the
# percentage-as-text Rate cells and the blank Status cells) two slots to the right,
# landing them correctly in H:R without needing to be re-written.
$ws.Columns("F").Insert()
$ws.Columns("H").Insert()

# Header row (row 1): rename D1, and set the two new headers F1 and H1.
# (The other headers, including the shifted ones, are already correct after the insert.)
$ws.Range("D1").Value = "Opening Balance"
$ws.Range("F1").Value = "Cost of asset sold"
$ws.Range("H1").Value = "Current Balance"

# Data rows 2-14: update Opening Balance (D), Cost of asset sold (F), Total (G),
# Current Balance (H), Rate (Pre MAB Assets) (I) and the depreciation figures (K:Q).
# Columns E, J ("Rate (Post MAB Assets)" text) and R ("Status", blank) are left untouched
# because they are unchanged by this edit and already shifted into place correctly.
# Row 2: Computer and Computer Equiments
$ws.Range("D2").Value = 2214149.3
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2214149.3
$ws.Range("H2").Value = 2214149.3
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 442829.86
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 442829.86
$ws.Range("P2").Value = 442829.86
$ws.Range("Q2").Value = 1771319.44

# Row 3: Office equipment
$ws.Range("D3").Value = 3623967
$ws.Range("F3").Value = 1686182
$ws.Range("G3").Value = 3623967
$ws.Range("H3").Value = 1937785
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 517709.57
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 517709.57
$ws.Range("P3").Value = 517709.57
$ws.Range("Q3").Value = 3106257.43

# Row 4: Furniture
$ws.Range("D4").Value = 3111004.93
$ws.Range("F4").Value = 1367021.79
$ws.Range("G4").Value = 3111004.93
$ws.Range("H4").Value = 1743983.14
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 311100.49
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 311100.49
$ws.Range("P4").Value = 311100.49
$ws.Range("Q4").Value = 2799904.44

# Row 5: Telecommunications
$ws.Range("D5").Value = 334523
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 334523
$ws.Range("H5").Value = 334523
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 47789
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 47789
$ws.Range("P5").Value = 47789
$ws.Range("Q5").Value = 286734

# Row 6: Motor Vehicles
$ws.Range("D6").Value = 21671646
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 21671646
$ws.Range("H6").Value = 21671646
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 2167164.6
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 2167164.6
$ws.Range("P6").Value = 2167164.6
$ws.Range("Q6").Value = 19504481.4

# Row 7: Civil Works
$ws.Range("D7").Value = 549298.2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 549298.2
$ws.Range("H7").Value = 549298.2
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 27464.91
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 27464.91
$ws.Range("P7").Value = 27464.91
$ws.Range("Q7").Value = 521833.29

# Row 8: Internet and PABX System 
$ws.Range("D8").Value = 718355
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 718355
$ws.Range("H8").Value = 718355
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 47890.33
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 47890.33
$ws.Range("P8").Value = 47890.33
$ws.Range("Q8").Value = 670464.67

# Row 9: Speaker System
$ws.Range("D9").Value = 66900
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 66900
$ws.Range("H9").Value = 66900
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 4460
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 4460
$ws.Range("P9").Value = 4460
$ws.Range("Q9").Value = 62440

# Row 10: Sanitary Fittings
$ws.Range("D10").Value = 607855
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 607855
$ws.Range("H10").Value = 607855
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 30392.75
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 30392.75
$ws.Range("P10").Value = 30392.75
$ws.Range("Q10").Value = 577462.25

# Row 11: Electrical Wire and Fittings
$ws.Range("D11").Value = 1402841.25
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 1402841.25
$ws.Range("H11").Value = 1402841.25
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 93522.75
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 93522.75
$ws.Range("P11").Value = 93522.75
$ws.Range("Q11").Value = 1309318.5

# Row 12: Paint Works
$ws.Range("D12").Value = 375968.75
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 375968.75
$ws.Range("H12").Value = 375968.75
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 37596.88
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 37596.88
$ws.Range("P12").Value = 37596.88
$ws.Range("Q12").Value = 338371.87

# Row 13: Interior Fit-out Works
$ws.Range("D13").Value = 4733279.250000001
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 4733279.25
$ws.Range("H13").Value = 4733279.25
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 236663.96
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 236663.96
$ws.Range("P13").Value = 236663.96
$ws.Range("Q13").Value = 4496615.29

# Row 14: Furniture (Leasehold)
$ws.Range("D14").Value = 4771077.5
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 4771077.5
$ws.Range("H14").Value = 4771077.5
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 318071.83
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 318071.83
$ws.Range("P14").Value = 318071.83
$ws.Range("Q14").Value = 4453005.67
